# Applies the scheduled market-data refresh to the Leve profit tables.
# For each affected leve row, currentAveragePrice(NQ/HQ), LevePrice(NQ/HQ)
# and LeveProfit(NQ/HQ) columns (H:N) are updated with refreshed values.
$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row 2
$ws.Cells.Item(2, 8).Value = 209.53847
$ws.Cells.Item(2, 9).Value = 231
$ws.Cells.Item(2, 10).Value = 200
$ws.Cells.Item(2, 11).Value = 231
$ws.Cells.Item(2, 12).Value = 200
$ws.Cells.Item(2, 13).Value = -118
$ws.Cells.Item(2, 14).Value = -426
# row 9
$ws.Cells.Item(9, 8).Value = 63200.5
$ws.Cells.Item(9, 10).Value = 1939.8
$ws.Cells.Item(9, 12).Value = 1939.8
$ws.Cells.Item(9, 14).Value = -2277.8
# row 18
$ws.Cells.Item(18, 8).Value = 5700.5
$ws.Cells.Item(18, 9).Value = 5700.5
$ws.Cells.Item(18, 10).Value = 0
$ws.Cells.Item(18, 11).Value = 5700.5
$ws.Cells.Item(18, 12).Value = 0
$ws.Cells.Item(18, 13).Value = -5416.5
$ws.Cells.Item(18, 14).ClearContents()
# row 62
$ws.Cells.Item(62, 8).Value = 52224252
$ws.Cells.Item(62, 9).Value = 65279264
$ws.Cells.Item(62, 11).Value = 65279264
$ws.Cells.Item(62, 13).Value = -65278640
# row 64
$ws.Cells.Item(64, 8).Value = 5450
$ws.Cells.Item(64, 9).Value = 4400
$ws.Cells.Item(64, 11).Value = 4400
$ws.Cells.Item(64, 13).Value = -4152
# row 65
$ws.Cells.Item(65, 8).Value = 52224252
$ws.Cells.Item(65, 9).Value = 65279264
$ws.Cells.Item(65, 11).Value = 326396320
$ws.Cells.Item(65, 13).Value = -326393200
# row 67
$ws.Cells.Item(67, 8).Value = 5450
$ws.Cells.Item(67, 9).Value = 4400
$ws.Cells.Item(67, 11).Value = 4400
$ws.Cells.Item(67, 13).Value = -3542
# row 101
$ws.Cells.Item(101, 8).Value = 3882.2778
$ws.Cells.Item(101, 9).Value = 1207.5834
$ws.Cells.Item(101, 11).Value = 3622.7502
$ws.Cells.Item(101, 13).Value = -2000.7502
# row 104
$ws.Cells.Item(104, 8).Value = 100
$ws.Cells.Item(104, 9).Value = 100
$ws.Cells.Item(104, 11).Value = 300
$ws.Cells.Item(104, 13).Value = 1447
# row 116
$ws.Cells.Item(116, 8).Value = 17721.957
$ws.Cells.Item(116, 10).Value = 18186.592
$ws.Cells.Item(116, 12).Value = 18186.592
$ws.Cells.Item(116, 14).Value = -25070.592
# row 132
$ws.Cells.Item(132, 8).Value = 4713.5713
$ws.Cells.Item(132, 9).Value = 4696
$ws.Cells.Item(132, 11).Value = 14088
$ws.Cells.Item(132, 13).Value = -11558
# row 138
$ws.Cells.Item(138, 8).Value = 3679.3513
$ws.Cells.Item(138, 10).Value = 4645.7144
$ws.Cells.Item(138, 12).Value = 13937.1432
$ws.Cells.Item(138, 14).Value = -24217.1432

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row 45
$ws.Cells.Item(45, 8).Value = 2223.6667
$ws.Cells.Item(45, 9).Value = 1532
$ws.Cells.Item(45, 10).Value = 2915.3333
$ws.Cells.Item(45, 11).Value = 1532
$ws.Cells.Item(45, 12).Value = 2915.3333
$ws.Cells.Item(45, 13).Value = -1155
$ws.Cells.Item(45, 14).Value = -3669.3333
# row 97
$ws.Cells.Item(97, 8).Value = 3275.3845
$ws.Cells.Item(97, 9).Value = 1494.375
$ws.Cells.Item(97, 10).Value = 6125
$ws.Cells.Item(97, 11).Value = 1494.375
$ws.Cells.Item(97, 12).Value = 6125
$ws.Cells.Item(97, 13).Value = -998.375
$ws.Cells.Item(97, 14).Value = -7117
# row 122
$ws.Cells.Item(122, 8).Value = 1730.5834
$ws.Cells.Item(122, 9).Value = 839.64703
$ws.Cells.Item(122, 11).Value = 2518.94109
$ws.Cells.Item(122, 13).Value = -68.9410899999998
# row 132
$ws.Cells.Item(132, 8).Value = 3538.8276
$ws.Cells.Item(132, 9).Value = 2438.5217
$ws.Cells.Item(132, 11).Value = 7315.5651
$ws.Cells.Item(132, 13).Value = -4785.5651

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# row 94
$ws.Cells.Item(94, 8).Value = 8352367
$ws.Cells.Item(94, 10).Value = 27761.75
$ws.Cells.Item(94, 12).Value = 27761.75
$ws.Cells.Item(94, 14).Value = -28663.75
# row 99
$ws.Cells.Item(99, 8).Value = 5954.5
$ws.Cells.Item(99, 9).Value = 7556.8335
$ws.Cells.Item(99, 11).Value = 7556.8335
$ws.Cells.Item(99, 13).Value = -6058.8335
# row 105
$ws.Cells.Item(105, 8).Value = 2545.4443
$ws.Cells.Item(105, 9).Value = 2139.923
$ws.Cells.Item(105, 10).Value = 3599.8
$ws.Cells.Item(105, 11).Value = 2139.923
$ws.Cells.Item(105, 12).Value = 3599.8
$ws.Cells.Item(105, 13).Value = -392.9229999999998
$ws.Cells.Item(105, 14).Value = -7093.8
# row 107
$ws.Cells.Item(107, 8).Value = 1744.9302
$ws.Cells.Item(107, 9).Value = 1820.25
$ws.Cells.Item(107, 11).Value = 1820.25
$ws.Cells.Item(107, 13).Value = 99.75
# row 134
$ws.Cells.Item(134, 8).Value = 11994.818
$ws.Cells.Item(134, 9).Value = 16157.182
$ws.Cells.Item(134, 11).Value = 48471.546
$ws.Cells.Item(134, 13).Value = -45936.546

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Cells.Item(31, 8).Value = 3531
$ws.Cells.Item(31, 9).Value = 1388.8334
$ws.Cells.Item(31, 11).Value = 1388.8334
$ws.Cells.Item(31, 13).Value = -1093.8334
# row 34
$ws.Cells.Item(34, 8).Value = 3531
$ws.Cells.Item(34, 9).Value = 1388.8334
$ws.Cells.Item(34, 11).Value = 1388.8334
$ws.Cells.Item(34, 13).Value = -1186.8334
# row 58
$ws.Cells.Item(58, 8).Value = 5481.8125
$ws.Cells.Item(58, 9).Value = 4566.222
$ws.Cells.Item(58, 10).Value = 6659
$ws.Cells.Item(58, 11).Value = 4566.222
$ws.Cells.Item(58, 12).Value = 6659
$ws.Cells.Item(58, 13).Value = -4363.222
$ws.Cells.Item(58, 14).Value = -7065
# row 99
$ws.Cells.Item(99, 8).Value = 2914.4644
$ws.Cells.Item(99, 9).Value = 2814.1177
$ws.Cells.Item(99, 10).Value = 3069.5454
$ws.Cells.Item(99, 11).Value = 2814.1177
$ws.Cells.Item(99, 12).Value = 3069.5454
$ws.Cells.Item(99, 13).Value = -1316.1177
$ws.Cells.Item(99, 14).Value = -6065.5454
# row 126
$ws.Cells.Item(126, 8).Value = 2914.4644
$ws.Cells.Item(126, 9).Value = 2814.1177
$ws.Cells.Item(126, 10).Value = 3069.5454
$ws.Cells.Item(126, 11).Value = 8442.3531
$ws.Cells.Item(126, 12).Value = 9208.636200000001
$ws.Cells.Item(126, 13).Value = -5972.3531
$ws.Cells.Item(126, 14).Value = -14148.6362
# row 132
$ws.Cells.Item(132, 8).Value = 2020
$ws.Cells.Item(132, 9).Value = 2045.4839
$ws.Cells.Item(132, 11).Value = 6136.4517
$ws.Cells.Item(132, 13).Value = -3606.4517
# row 134
$ws.Cells.Item(134, 8).Value = 5144.069
$ws.Cells.Item(134, 9).Value = 4570.587
$ws.Cells.Item(134, 11).Value = 13711.761
$ws.Cells.Item(134, 13).Value = -11176.761
# row 136
$ws.Cells.Item(136, 8).Value = 5481.8125
$ws.Cells.Item(136, 9).Value = 4566.222
$ws.Cells.Item(136, 10).Value = 6659
$ws.Cells.Item(136, 11).Value = 13698.666
$ws.Cells.Item(136, 12).Value = 19977
$ws.Cells.Item(136, 13).Value = -11148.666
$ws.Cells.Item(136, 14).Value = -25077

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# row 129
$ws.Cells.Item(129, 8).Value = 13891692
$ws.Cells.Item(129, 10).Value = 20836746
$ws.Cells.Item(129, 12).Value = 62510238
$ws.Cells.Item(129, 14).Value = -62520238
# row 131
$ws.Cells.Item(131, 8).Value = 3592.53
$ws.Cells.Item(131, 9).Value = 3233.3333
$ws.Cells.Item(131, 10).Value = 3615.4575
$ws.Cells.Item(131, 11).Value = 9699.999899999999
$ws.Cells.Item(131, 12).Value = 10846.3725
$ws.Cells.Item(131, 13).Value = -4659.999899999999
$ws.Cells.Item(131, 14).Value = -20926.3725
# row 134
$ws.Cells.Item(134, 8).Value = 11801.823
$ws.Cells.Item(134, 9).Value = 5966.4546
$ws.Cells.Item(134, 11).Value = 17899.3638
$ws.Cells.Item(134, 13).Value = -12829.3638

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# row 70
$ws.Cells.Item(70, 8).Value = 6263.433
$ws.Cells.Item(70, 9).Value = 6108.625
$ws.Cells.Item(70, 10).Value = 6882.6665
$ws.Cells.Item(70, 11).Value = 6108.625
$ws.Cells.Item(70, 12).Value = 6882.6665
$ws.Cells.Item(70, 13).Value = -5838.625
$ws.Cells.Item(70, 14).Value = -7422.6665
# row 73
$ws.Cells.Item(73, 8).Value = 6263.433
$ws.Cells.Item(73, 9).Value = 6108.625
$ws.Cells.Item(73, 10).Value = 6882.6665
$ws.Cells.Item(73, 11).Value = 6108.625
$ws.Cells.Item(73, 12).Value = 6882.6665
$ws.Cells.Item(73, 13).Value = -5172.625
$ws.Cells.Item(73, 14).Value = -8754.666499999999
# row 80
$ws.Cells.Item(80, 8).Value = 2069.889
$ws.Cells.Item(80, 9).Value = 2168.6
$ws.Cells.Item(80, 10).Value = 1946.5
$ws.Cells.Item(80, 11).Value = 2168.6
$ws.Cells.Item(80, 12).Value = 1946.5
$ws.Cells.Item(80, 13).Value = -1170.6
$ws.Cells.Item(80, 14).Value = -3942.5
# row 83
$ws.Cells.Item(83, 8).Value = 2069.889
$ws.Cells.Item(83, 9).Value = 2168.6
$ws.Cells.Item(83, 10).Value = 1946.5
$ws.Cells.Item(83, 11).Value = 10843
$ws.Cells.Item(83, 12).Value = 9732.5
$ws.Cells.Item(83, 13).Value = -5851
$ws.Cells.Item(83, 14).Value = -19716.5
# row 132
$ws.Cells.Item(132, 8).Value = 5929.2
$ws.Cells.Item(132, 9).Value = 3859.2
$ws.Cells.Item(132, 10).Value = 7999.2
$ws.Cells.Item(132, 11).Value = 11577.6
$ws.Cells.Item(132, 12).Value = 23997.6
$ws.Cells.Item(132, 13).Value = -9047.599999999999
$ws.Cells.Item(132, 14).Value = -29057.6

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# row 22
$ws.Cells.Item(22, 8).Value = 1197.909
$ws.Cells.Item(22, 10).Value = 1197.909
$ws.Cells.Item(22, 12).Value = 1197.909
$ws.Cells.Item(22, 14).Value = -1787.909
# row 27
$ws.Cells.Item(27, 8).Value = 1197.909
$ws.Cells.Item(27, 10).Value = 1197.909
$ws.Cells.Item(27, 12).Value = 1197.909
$ws.Cells.Item(27, 14).Value = -1411.909
# row 40
$ws.Cells.Item(40, 8).Value = 2256.4167
$ws.Cells.Item(40, 9).Value = 2067.9
$ws.Cells.Item(40, 11).Value = 2067.9
$ws.Cells.Item(40, 13).Value = -1931.9
# row 46
$ws.Cells.Item(46, 8).Value = 2166.1304
$ws.Cells.Item(46, 10).Value = 2501.3333
$ws.Cells.Item(46, 12).Value = 2501.3333
$ws.Cells.Item(46, 14).Value = -2877.3333
# row 93
$ws.Cells.Item(93, 8).Value = 1485.7142
$ws.Cells.Item(93, 9).Value = 1172.9
$ws.Cells.Item(93, 11).Value = 1172.9
$ws.Cells.Item(93, 13).Value = 75.09999999999991
# row 122
$ws.Cells.Item(122, 8).Value = 6020.433
$ws.Cells.Item(122, 9).Value = 4552.857
$ws.Cells.Item(122, 10).Value = 6467.087
$ws.Cells.Item(122, 11).Value = 13658.571
$ws.Cells.Item(122, 12).Value = 19401.261
$ws.Cells.Item(122, 13).Value = -11208.571
$ws.Cells.Item(122, 14).Value = -24301.261
# row 136
$ws.Cells.Item(136, 8).Value = 1783.4445
$ws.Cells.Item(136, 9).Value = 1221.4
$ws.Cells.Item(136, 11).Value = 3664.2
$ws.Cells.Item(136, 13).Value = -1114.2
# row 138
$ws.Cells.Item(138, 8).Value = 71884
$ws.Cells.Item(138, 10).Value = 70428
$ws.Cells.Item(138, 12).Value = 70428
$ws.Cells.Item(138, 14).Value = -80708

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# row 113
$ws.Cells.Item(113, 8).Value = 346.1111
$ws.Cells.Item(113, 9).Value = 349.5
$ws.Cells.Item(113, 10).Value = 339.33334
$ws.Cells.Item(113, 11).Value = 1048.5
$ws.Cells.Item(113, 12).Value = 1018.00002
$ws.Cells.Item(113, 13).Value = 1121.5
$ws.Cells.Item(113, 14).Value = -5358.00002
# row 122
$ws.Cells.Item(122, 8).Value = 242454.86
$ws.Cells.Item(122, 9).Value = 337054.66
$ws.Cells.Item(122, 11).Value = 1011163.98
$ws.Cells.Item(122, 13).Value = -1008713.98
# row 126
$ws.Cells.Item(126, 8).Value = 2474.5
$ws.Cells.Item(126, 9).Value = 2659.2
$ws.Cells.Item(126, 10).Value = 2166.6667
$ws.Cells.Item(126, 11).Value = 7977.599999999999
$ws.Cells.Item(126, 12).Value = 6500.000100000001
$ws.Cells.Item(126, 13).Value = -5507.599999999999
$ws.Cells.Item(126, 14).Value = -11440.0001
# row 136
$ws.Cells.Item(136, 8).Value = 9870.931
$ws.Cells.Item(136, 9).Value = 9866.464
$ws.Cells.Item(136, 11).Value = 29599.392
$ws.Cells.Item(136, 13).Value = -27049.392
